$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The "Förändrad" (last-changed) date for every existing data row (2..423)
#    moved from 2023-09-06 (45175) to 2023-09-08 (45177).
$ws.Range("C2:C423").Value = 45177

# 2) Row 423 picks up an explicit row height, matching the customHeight flag
#    already present on every other data row.
$ws.Rows.Item(423).RowHeight = 15

# 3) A brand-new case was appended as row 424.
$ws.Range("A424").Value = "A 41440-2023"

$ws.Range("B424").Value = 45175
$ws.Range("B424").NumberFormat = "YYYY-MM-DD"

$ws.Range("C424").Value = 45177
$ws.Range("C424").NumberFormat = "YYYY-MM-DD"

$ws.Range("D424").Value = "SKÅNE LÄN"
$ws.Range("E424").Value = "OSBY"

$ws.Range("G424").Value = 0.6
$ws.Range("H424").Value = 0
$ws.Range("I424").Value = 0
$ws.Range("J424").Value = 0
$ws.Range("K424").Value = 0
$ws.Range("L424").Value = 0
$ws.Range("M424").Value = 0
$ws.Range("N424").Value = 0
$ws.Range("O424").Value = 0
$ws.Range("P424").Value = 0
$ws.Range("Q424").Value = 0

# R424 stays empty (no species found) but keeps the wrap-text style used by
# every other "no species" row in column R.
$ws.Range("R424").WrapText = $true
